$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The project schedule had two milestone rows removed ("Hito 3: Desarrollo de
# la Base de datos" and "Hito 5: Pruebas de calidad QA"), the remaining
# milestones were renumbered/retitled, and the final milestone picked up the
# completion date that used to sit on the (now removed) trailing blank cell.
# ---------------------------------------------------------------------------

# Delete the higher-numbered row first so the lower row number used below
# still points at the right row.
$ws.Rows("33:33").Delete()
$ws.Rows("27:27").Delete()

# Rename the milestone that used to be "Hito 4: Desarrollo Front-End y
# Back-End" (now at row 29) to its corrected title.
$ws.Range("B29").Value = "Hito 3: Desarrollo y Pruebas"

# Rename the final milestone (now at row 34) and give it the completion date
# that the preceding activity row already carries.
$ws.Range("B34").Value = "Hito 4: Despliegue"
$ws.Range("C34").Value = $ws.Range("G33").Value()

# Re-apply the AutoFilter over the new (shorter) table range.
$ws.AutoFilterMode = $false
$ws.Range("B10:H34").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range.
$name = $wb.Names.Item(1)
$name.RefersTo = "='Hoja 1'!`$B`$10:`$H`$34"
